$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: 330uF Tantalum Cap
$url1 = "https://www.mouser.co.uk/ProductDetail/Panasonic/4TPF330MFL?qs=sGAEpiMZZMsh%252B1woXyUXj9v%2FLaJorgg0K34QuPhlki8%3D"
$ws.Range("D13").Value = $url1
$ws.Hyperlinks.Add($ws.Range("D13"), $url1)
$ws.Range("D13").Style = "Hyperlink"

$ws.Range("C13").Value = "330uF Tantalum Cap"

# Row 14: 47uF Ceramic Cap
$ws.Range("C14").Value = "47uF Ceramic Cap"

$url2 = "https://www.mouser.co.uk/ProductDetail/TDK/CGA9N1X7R1V476M230KC?qs=sGAEpiMZZMsh%252B1woXyUXj2GDJWaunJJjOXJcg%252BYyWqY%3D"
$ws.Range("D14").Value = $url2
$ws.Hyperlinks.Add($ws.Range("D14"), $url2)
$ws.Range("D14").Style = "Hyperlink"

$ws.Range("E14").Value = "5750(metric)"

$ws.Range("E13").Value = 7343
$ws.Range("E13").HorizontalAlignment = -4131

$ws.Range("F23").Select()
